# New crime data collected - weekly CompStat report refresh
# Updates the "Volume/Number" + "Week covering" header text, the 022 precinct
# weekly crime table (rows 15,16,18,19,21,24,25) with the new week's figures,
# and the auto-fit width of column E that results from the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: bump the report volume/number and the covered week's dates.
# Both cells are single shared-string cells made of several same-formatted
# rich-text runs, so writing the whole merged string reproduces them exactly.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/13/2023  Through  11/19/2023"

# ---------------------------------------------------------------------------
# Helper pattern used below for cells that change numeric <-> text "category"
# (and therefore need a different cell style): copy the *format* only from a
# same-shaped cell that already carries the destination style, then write the
# new value. PasteSpecial(-4122) = xlPasteFormats, value-preserving.
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

# Row 15 ---------------------------------------------------------------
$ws.Range("N15").Value = -80

# Row 16 ---------------------------------------------------------------
$ws.Range("F16").Copy()
$ws.Range("C16").PasteSpecial($xlPasteFormats)
$ws.Range("C16").Value = 1

$ws.Range("G16").Value = 3
$ws.Range("H16").Value = -66.666666666666
$ws.Range("I16").Value = 19
$ws.Range("K16").Value = -32.142857142857
$ws.Range("L16").Value = 35.714285714285
$ws.Range("M16").Value = -13.636363636363
$ws.Range("N16").Value = -89.673913043478

# Row 18 ---------------------------------------------------------------
$ws.Range("N18").Value = -84.848484848484

# Row 19 ---------------------------------------------------------------
$ws.Range("F16").Copy()
$ws.Range("C19").PasteSpecial($xlPasteFormats)
$ws.Range("C19").Value = 1

# D19 becomes the literal text "0" (numeric-looking) - force text storage via
# a Text number format before assignment, then restore the normal "General"
# text style used elsewhere in the table (copied from C26, which keeps it).
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0"
$ws.Range("C26").Copy()
$ws.Range("D19").PasteSpecial($xlPasteFormats)

$ws.Range("C26").Copy()
$ws.Range("E19").PasteSpecial($xlPasteFormats)
$ws.Range("E19").Value = "***.*"

$ws.Range("F19").Value = 2
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 46
$ws.Range("K19").Value = 84
$ws.Range("L19").Value = 100
$ws.Range("M19").Value = -24.590163934426
$ws.Range("N19").Value = -72.781065088757

# Row 21 (bold TOTAL row) ------------------------------------------------
$ws.Range("C43").Copy()
$ws.Range("C21").PasteSpecial($xlPasteFormats)
$ws.Range("C21").Value = 2

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0"
$ws.Range("A21").Copy()
$ws.Range("D21").PasteSpecial($xlPasteFormats)

$ws.Range("A21").Copy()
$ws.Range("E21").PasteSpecial($xlPasteFormats)
$ws.Range("E21").Value = "***.*"

$ws.Range("F21").Value = 3
$ws.Range("G21").Value = 6
$ws.Range("H21").Value = -50
$ws.Range("I21").Value = 80
$ws.Range("K21").Value = 11.111111111111
$ws.Range("L21").Value = 50.943396226415
$ws.Range("M21").Value = -17.525773195876
$ws.Range("N21").Value = -81.981981981982

# Row 24 ---------------------------------------------------------------
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "0"
$ws.Range("C26").Copy()
$ws.Range("C24").PasteSpecial($xlPasteFormats)

$ws.Range("M24").Value = -60.416666666666

# Row 25 ---------------------------------------------------------------
$ws.Range("F16").Copy()
$ws.Range("C25").PasteSpecial($xlPasteFormats)
$ws.Range("C25").Value = 2

$ws.Range("F16").Copy()
$ws.Range("D25").PasteSpecial($xlPasteFormats)
$ws.Range("D25").Value = 1

$ws.Range("N15").Copy()
$ws.Range("E25").PasteSpecial($xlPasteFormats)
$ws.Range("E25").Value = 100

$ws.Range("F25").Value = 3
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 50
$ws.Range("I25").Value = 47
$ws.Range("J25").Value = 29
$ws.Range("K25").Value = 62.068965517241
$ws.Range("L25").Value = 104.347826086957
$ws.Range("M25").Value = 161.111111111111

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Column E auto-fit width shrinks now that E19/E21 hold the shorter "***.*"
# text instead of "-100". Reproduce the recalculated best-fit width (closest
# reachable value to the canonical 7.433768 via the ColumnWidth property).
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 6.714285714285714
